$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revert the accented names that were previously stripped of accents.
$ws.Range("A2").Value = "Ignacio Fernández Fernández"
$ws.Range("A3").Value = "Naucé López González"

# Restore the selection to C2 (as in the original/target file).
$ws.Range("C2").Select()
